# Login.xlsx - add Chat Connection / Location Filter / Export Deal / Negotiate
# Deal test-case rows (66-79) to the "Login" sheet, and move the active
# selection down to the newly-added block (mirrors the author's commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# --- Rows 66-69: "Chat Connection" test cases -------------------------------
# Filled column-by-column (A66:A69 first, then B66:B69) the way the original
# author's sheet was edited.
$ws.Cells.Item(66, 1).Value = "ChatConnection_TC001"
$ws.Cells.Item(67, 1).Value = "ChatConnection_TC001(2)"
$ws.Cells.Item(68, 1).Value = "ChatConnection_TC002"
$ws.Cells.Item(69, 1).Value = "ChatConnection_TC002(2)"

$ws.Cells.Item(66, 2).Value = "rogerdeals21+rick@gmail.com"
$ws.Cells.Item(67, 2).Value = "rogerdeals21+john@gmail.com"
$ws.Cells.Item(68, 2).Value = "rogerdeals21+matt@gmail.com"
$ws.Cells.Item(69, 2).Value = "rogerdeals21+zeb@gmail.com"

for ($r = 66; $r -le 69; $r++) {
    $ws.Cells.Item($r, 3).Value = "arewethere?"
    $ws.Cells.Item($r, 4).Value = "Login successful"
}

# --- Rows 70-79: Location Filter / Export Deal / Negotiate Deal test cases -
# Columns: A = Automation Test ID, B = UserName, C = Password, D = Expected Result
$rows = @(
    @("LocFilter_TC001",        "rogerdeals21+stan@gmail.com"),
    @("LocFilter_TC002",        "rogerdeals21+rick@gmail.com"),
    @("LocFilter_TC003",        "rogerdeals21+john@gmail.com"),
    @("ExportDeal_TC001",       "rogerdeals21+stan@gmail.com"),
    @("ExportDeal_TC002",       "rogerdeals21+rick@gmail.com"),
    @("NegotiateDeal_TC001",    "rogerdeals21+stan@gmail.com"),
    @("NegotiateDeal_TC002",    "rogerdeals21+rick@gmail.com"),
    @("NegotiateDeal_TC003",    "rogerdeals21+john@gmail.com"),
    @("NegotiateDeal_TC004",    "rogerdeals21+stan@gmail.com"),
    @("NegotiateDeal_TC004(2)", "rogerdeals21+john@gmail.com")
)

$startRow = 70
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = "arewethere?"
    $ws.Cells.Item($r, 4).Value = "Login successful"

    $ws.Cells.Item($r, 1).VerticalAlignment = -4108
}

# --- Move the selection / viewport to the newly added block ----------------
[void]$ws.Range("A70:D79").Select()

Write-Output "Added rows 66-79 (Chat Connection / Loc Filter / Export Deal / Negotiate Deal test cases)"
